$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 2 de Mayo de 2020 a las 23:41'
$ws.Cells.Item(4, 2).Value = 62073
$ws.Cells.Item(4, 3).Value = 37530
$ws.Cells.Item(4, 4).Value = 16251
$ws.Cells.Item(4, 5).Value = 8292
$ws.Cells.Item(5, 2).Value = 49850
$ws.Cells.Item(5, 3).Value = 19913
$ws.Cells.Item(5, 4).Value = 24800
$ws.Cells.Item(5, 5).Value = 5137
$ws.Cells.Item(6, 2).Value = 17165
$ws.Cells.Item(6, 3).Value = 6813
$ws.Cells.Item(6, 4).Value = 8564
$ws.Cells.Item(6, 5).Value = 1788
$ws.Cells.Item(7, 2).Value = 15967
$ws.Cells.Item(7, 3).Value = 5702
$ws.Cells.Item(7, 4).Value = 7731
$ws.Cells.Item(7, 5).Value = 2534
$ws.Cells.Item(9, 2).Value = 12207
$ws.Cells.Item(9, 3).Value = 6821
$ws.Cells.Item(9, 4).Value = 4133
$ws.Cells.Item(9, 5).Value = 1253
$ws.Cells.Item(10, 2).Value = 8902
$ws.Cells.Item(10, 3).Value = 5981
$ws.Cells.Item(10, 4).Value = 2364
$ws.Cells.Item(10, 5).Value = 557
$ws.Cells.Item(13, 1).Value = 'Aragon'
$ws.Cells.Item(13, 2).Value = 5136
$ws.Cells.Item(13, 3).Value = 2474
$ws.Cells.Item(13, 4).Value = 1901
$ws.Cells.Item(13, 5).Value = 761
$ws.Cells.Item(14, 1).Value = 'Valencia/Valencia'
$ws.Cells.Item(14, 2).Value = 5131
$ws.Cells.Item(14, 3).Value = 2194
$ws.Cells.Item(14, 4).Value = 2583
$ws.Cells.Item(14, 5).Value = 515
$ws.Cells.Item(15, 2).Value = 4888
$ws.Cells.Item(15, 3).Value = 2281
$ws.Cells.Item(15, 4).Value = 2147
$ws.Cells.Item(15, 5).Value = 460
$ws.Cells.Item(16, 1).Value = 'La Rioja'
$ws.Cells.Item(16, 2).Value = 3947
$ws.Cells.Item(16, 3).Value = 2279
$ws.Cells.Item(16, 4).Value = 1393
$ws.Cells.Item(16, 5).Value = 334
$ws.Cells.Item(17, 1).Value = 'Toledo'
$ws.Cells.Item(17, 2).Value = 3938
$ws.Cells.Item(17, 3).Value = 4178
$ws.Cells.Item(17, 4).Value = 10597
$ws.Cells.Item(17, 5).Value = 504
$ws.Cells.Item(23, 2).Value = 2824
$ws.Cells.Item(23, 3).Value = 2017
$ws.Cells.Item(23, 4).Value = 349
$ws.Cells.Item(23, 5).Value = 458
$ws.Cells.Item(30, 2).Value = 2303
$ws.Cells.Item(30, 3).Value = 890
$ws.Cells.Item(30, 4).Value = 1134
$ws.Cells.Item(30, 5).Value = 279
$ws.Cells.Item(32, 2).Value = 2212
$ws.Cells.Item(32, 3).Value = 1149
$ws.Cells.Item(32, 4).Value = 923
$ws.Cells.Item(32, 5).Value = 140
$ws.Cells.Item(33, 2).Value = 2201
$ws.Cells.Item(33, 3).Value = 1551
$ws.Cells.Item(33, 4).Value = 457
$ws.Cells.Item(33, 5).Value = 193
$ws.Cells.Item(38, 2).Value = 1491
$ws.Cells.Item(38, 3).Value = 1251
$ws.Cells.Item(38, 4).Value = 108
